$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_StatusQry")

# Update the "Last updated" timestamp banner in A1
$ws.Range("A1").Value = "Last updated: 2025-07-15 06:10:51"

# Row 9 (PO 4516351202_ARD): CommittedNotShip and ForPackOrders swap, HasUnshippedCommitted flips to -1
$ws.Range("C9").Value = 6
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = -1

# Row 20 (PO 4516351202_TIPI): CommittedNotShip and ForPackOrders swap, HasUnshippedCommitted flips to -1
$ws.Range("C20").Value = 13
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = -1
